$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing used range (A1:B4) and write the single new value
$ws.Range("A1:B4").Clear()
$ws.Range("A1").Value = "G607 DAN"

# Update the selection to match the target (F6)
$ws.Range("F6").Select()
